$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update topic text for the midterm weeks (row 6 and row 12 in sheetData,
# which are the C6 and C12 cells on the worksheet) to the new unified text.
$ws.Range("C6").Value = "Midterm on Wednesday"
$ws.Range("C12").Value = "Midterm on Wednesday"

# Update the active selection to match the saved workbook view state.
$ws.Range("C18").Select()
